# Update the Lorenz curve values on the Sim_5yr and Sim_10yr sheets.
$wb = $excel.ActiveWorkbook

$wsSim5yr  = $wb.Worksheets.Item("Sim_5yr")
$wsSim10yr = $wb.Worksheets.Item("Sim_10yr")

# --- Sim_5yr (sheet1): rows 2-10, columns B:E -----------------------------
$sim5yrValues = @{
    2 = @(0.001582864817699225, 0.03245530933421766, 0.09779696798524537, 0.2372381476731427)
    3 = @(-0.0007459087716984048, 0.03672921350623698, 0.1168834235655709, 0.2813652226324158)
    4 = @(0.0006011783247955002, 0.04247450218029081, 0.1307407285424236, 0.3104676561855404)
    5 = @(0.009372101279834375, 0.05798743441832001, 0.1542437689779933, 0.3435184102535908)
    6 = @(0.01824142196500853, 0.07246016959575424, 0.1751103290949647, 0.3684927036409594)
    7 = @(0.02323333857990465, 0.08145261189390031, 0.1877971626268973, 0.3806085007891786)
    8 = @(0.02408598800423544, 0.08205584913811813, 0.1875675541666817, 0.3796722118029505)
    9 = @(0.02325050282488948, 0.07944574816628962, 0.1829000877260878, 0.3725327601800408)
    10 = @(0.02208711002900622, 0.07632242826641517, 0.1769482614720676, 0.3638438421333888)
}

foreach ($row in $sim5yrValues.Keys) {
    $vals = $sim5yrValues[$row]
    $wsSim5yr.Cells.Item($row, 2).Value = $vals[0]
    $wsSim5yr.Cells.Item($row, 3).Value = $vals[1]
    $wsSim5yr.Cells.Item($row, 4).Value = $vals[2]
    $wsSim5yr.Cells.Item($row, 5).Value = $vals[3]
}

# --- Sim_10yr (sheet2): rows 2-6, columns B:E -----------------------------
$sim10yrValues = @{
    2 = @(0.001582864817699225, 0.03245530933421766, 0.09779696798524537, 0.2372381476731427)
    3 = @(-0.0002659506110192058, 0.03902747974518321, 0.1226509462066624, 0.2946909076209057)
    4 = @(0.01342847877888379, 0.06391894172670295, 0.1625751596521325, 0.3536324974972478)
    5 = @(0.02331079150650988, 0.08086837257482735, 0.1859033580011066, 0.37761309490064)
    6 = @(0.02265916598987138, 0.07791524785385777, 0.1799996033585577, 0.3683604986178114)
}

foreach ($row in $sim10yrValues.Keys) {
    $vals = $sim10yrValues[$row]
    $wsSim10yr.Cells.Item($row, 2).Value = $vals[0]
    $wsSim10yr.Cells.Item($row, 3).Value = $vals[1]
    $wsSim10yr.Cells.Item($row, 4).Value = $vals[2]
    $wsSim10yr.Cells.Item($row, 5).Value = $vals[3]
}

$wb.Save()
